$d = $word.ActiveDocument

# --- 1. Title: "while and do while" -> "while Loop" -----------------------
# The title paragraph starts life as 7 separate runs:
#   "while" " " "and" " " "do" " " "while"
# The target keeps the first two runs ("while" and " ") untouched, turns
# the "and" run into "Loop", and drops the trailing " do while" runs.
# Editing this engine's Range merges any run that is only partially
# touched into its left neighbour, but deleting (or inserting at) an
# *exact* run boundary leaves the surrounding runs alone - so do the
# trim in two exact-boundary deletes, then append the replacement text
# as a brand new run at the (now shorter) end of the paragraph.

$titlePar = $d.Paragraphs.Item(1).Range

$rAnd = $d.Range($titlePar.Start, $titlePar.End)
$rAnd.Find.Execute("and", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Remove the trailing " do while" (everything after "and", up to the end
# of the paragraph's text, i.e. not including the paragraph mark).
$tailStart = $rAnd.End
$tailEnd = $titlePar.End - 1
$rTail = $d.Range($tailStart, $tailEnd)
$rTail.Delete()

# Remove the now-isolated "and" run.
$rAndExact = $d.Range($rAnd.Start, $rAnd.End)
$rAndExact.Delete()

# Append "Loop" as its own run right after "while ".
$insPos = $d.Paragraphs.Item(1).Range.End - 1
$rIns = $d.Range($insPos, $insPos)
$rIns.InsertAfter("Loop")

# --- 2. Date stamp ----------------------------------------------------------
$d.Content.Find.Execute("October  27, 2021 (03:18:52 PM)", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "October  29, 2021 (07:12:55 PM)", 2)
